# Update "想去人数" (F column) values on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets, per the generated data refresh.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value  = 445
$ws1.Range("F7").Value  = 2216
$ws1.Range("F11").Value = 5005
$ws1.Range("F16").Value = 32
$ws1.Range("F21").Value = 3971
$ws1.Range("F22").Value = 722
$ws1.Range("F23").Value = 698
$ws1.Range("F34").Value = 1014
$ws1.Range("F36").Value = 2571
$ws1.Range("F38").Value = 25

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 445
$ws4.Range("F7").Value  = 2216
$ws4.Range("F11").Value = 5005
$ws4.Range("F16").Value = 32
$ws4.Range("F21").Value = 3971
$ws4.Range("F22").Value = 722
$ws4.Range("F23").Value = 698
$ws4.Range("F35").Value = 1014
$ws4.Range("F37").Value = 2571
$ws4.Range("F39").Value = 25
